$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current values for columns D, J, K, L, M, P across rows 2-14
$rows = 2..14
$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
}

# Mapping: row r receives the old (snapshot) values of row $perm[r]
$perm = @{
    2 = 13
    3 = 11
    4 = 6
    5 = 9
    6 = 5
    7 = 14
    8 = 8
    9 = 10
    10 = 3
    11 = 2
    12 = 4
    13 = 7
    14 = 12
}

foreach ($r in $rows) {
    $src = $perm[$r]
    $vals = $snapshot[$src]
    $ws.Cells.Item($r, 4).Value2 = $vals.D
    $ws.Cells.Item($r, 10).Value2 = $vals.J
    $ws.Cells.Item($r, 11).Value2 = $vals.K
    $ws.Cells.Item($r, 12).Value2 = $vals.L
    $ws.Cells.Item($r, 13).Value2 = $vals.M
    $ws.Cells.Item($r, 16).Value2 = $vals.P
}
